$d = $word.ActiveDocument
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphBefore()
$p6.Range.InsertParagraphBefore()
$deckp = $d.Paragraphs.Item(6)
$linkp = $d.Paragraphs.Item(7)
$linkp.Style = "Normal"

$deckp.Range.Text = "Deck"

$r = $d.Range($linkp.Range.Start, $linkp.Range.Start)
$addr = "https://github.com/SteveLasker/Presentations/tree/master/DotNetConf2016"
$d.Hyperlinks.Add($r, $addr, [Type]::Missing, [Type]::Missing, $addr) | Out-Null

$linkp2 = $d.Paragraphs.Item(7)
Write-Host ("linkp2 range: " + $linkp2.Range.Start + " " + $linkp2.Range.End)
$tailR = $d.Range($linkp2.Range.End - 1, $linkp2.Range.End - 1)
Write-Host ("tailR: " + $tailR.Start + " " + $tailR.End)
$tailR.InsertAfter(" ")

for ($i = 1; $i -le 10; $i++) {
    Write-Host ("[$i] '" + $d.Paragraphs.Item($i).Range.Text.Replace("`r","\r") + "'")
}
